# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed "K" column (column G) values, row -> new value
$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 1
    13 = 3
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 0
    22 = 1
    23 = 1
    24 = 1
    25 = 2
    26 = 1
    27 = 1
    28 = 1
    29 = 1
    30 = 0
    31 = 1
    32 = 3
    33 = 2
    34 = 2
    35 = 0
    36 = 3
    37 = 0
    38 = 2
    39 = 2
    40 = 1
    41 = 0
    42 = 1
    43 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
